$d = $word.ActiveDocument

$d.Content.Find.Execute("386÷2=193, 0", $true, $false, $false, $false, $false, $true, 1, $false, "691÷5=138, 1", 2) | Out-Null
$d.Content.Find.Execute("735÷3=245, 0", $true, $false, $false, $false, $false, $true, 1, $false, "463÷3=154, 1", 2) | Out-Null
$d.Content.Find.Execute("681÷9=75, 6", $true, $false, $false, $false, $false, $true, 1, $false, "226÷7=32, 2", 2) | Out-Null
$d.Content.Find.Execute("804÷6=134, 0", $true, $false, $false, $false, $false, $true, 1, $false, "814÷6=135, 4", 2) | Out-Null
$d.Content.Find.Execute("998÷9=110, 8", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=256, 1", 2) | Out-Null
$d.Content.Find.Execute("190÷4=47, 2", $true, $false, $false, $false, $false, $true, 1, $false, "760÷2=380, 0", 2) | Out-Null
$d.Content.Find.Execute("759÷6=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "652÷2=326, 0", 2) | Out-Null
$d.Content.Find.Execute("116÷8=14, 4", $true, $false, $false, $false, $false, $true, 1, $false, "778÷6=129, 4", 2) | Out-Null
$d.Content.Find.Execute("148÷4=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "658÷6=109, 4", 2) | Out-Null
$d.Content.Find.Execute("869÷9=96, 5", $true, $false, $false, $false, $false, $true, 1, $false, "951÷3=317, 0", 2) | Out-Null
$d.Content.Find.Execute("289÷4=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "139÷4=34, 3", 2) | Out-Null
$d.Content.Find.Execute("847÷5=169, 2", $true, $false, $false, $false, $false, $true, 1, $false, "256÷5=51, 1", 2) | Out-Null
$d.Content.Find.Execute("170÷9=18, 8", $true, $false, $false, $false, $false, $true, 1, $false, "665÷6=110, 5", 2) | Out-Null
$d.Content.Find.Execute("102÷2=51, 0", $true, $false, $false, $false, $false, $true, 1, $false, "475÷4=118, 3", 2) | Out-Null
$d.Content.Find.Execute("618÷3=206, 0", $true, $false, $false, $false, $false, $true, 1, $false, "189÷6=31, 3", 2) | Out-Null
$d.Content.Find.Execute("338÷9=37, 5", $true, $false, $false, $false, $false, $true, 1, $false, "770÷6=128, 2", 2) | Out-Null
$d.Content.Find.Execute("628÷9=69, 7", $true, $false, $false, $false, $false, $true, 1, $false, "780÷6=130, 0", 2) | Out-Null
$d.Content.Find.Execute("568÷2=284, 0", $true, $false, $false, $false, $false, $true, 1, $false, "735÷6=122, 3", 2) | Out-Null
$d.Content.Find.Execute("930÷5=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "342÷7=48, 6", 2) | Out-Null
$d.Content.Find.Execute("954÷8=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "226÷6=37, 4", 2) | Out-Null
$d.Content.Find.Execute("199÷8=24, 7", $true, $false, $false, $false, $false, $true, 1, $false, "477÷2=238, 1", 2) | Out-Null
$d.Content.Find.Execute("767÷9=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "314÷7=44, 6", 2) | Out-Null
$d.Content.Find.Execute("555÷7=79, 2", $true, $false, $false, $false, $false, $true, 1, $false, "398÷7=56, 6", 2) | Out-Null
$d.Content.Find.Execute("759÷7=108, 3", $true, $false, $false, $false, $false, $true, 1, $false, "476÷8=59, 4", 2) | Out-Null
$d.Content.Find.Execute("124÷3=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "291÷6=48, 3", 2) | Out-Null
